$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bText = @(
  'Description',
  'Unions representing workers at Turner   Newall say they are ''disappointed'' after talks with stricken parent firm Federal Mogul.',
  'SPACE.com - TORONTO, Canada -- A second\team of rocketeers competing for the  #36',
  'AP - A company founded by a chemistry researcher at the University of Louisville won a grant to develop a method of producing better peptides, which are short chains of amino acids, the building blocks of proteins.',
  'AP - It''s barely dawn when Mike Fitzpatrick starts his shift with a blur of colorful maps, figures and endless charts, but already he knows what the day will bring. Lightning will strike in places he expects. Winds will pick up, moist places will dry and flames will roar.',
  'AP - Southern California''s smog-fighting agency went after emissions of the bovine variety Friday, adopting the nation''s first rules to reduce air pollution from dairy cow manure.',
  'The British Department for Education and Skills (DfES) recently launched a "Music Manifesto" campaign, with the ostensible intention of educating the next generation of British musicians. Unfortunately, they also teamed up with the music industry (EMI, and various artists) to make this popular. EMI has apparently negotiated their end well, so that children in our schools will now be indoctrinated about the illegality of downloading music.The ignorance and audacity of this got to me a little, so I wrote an open letter to the DfES about it. Unfortunately, it''s pedantic, as I suppose you have to be when writing to goverment representatives. But I hope you find it useful, and perhaps feel inspired to do something similar, if or when the same thing has happened in your area.',
  '\\"Sven Jaschan, self-confessed author of the Netsky and Sasser viruses, is\responsible for 70 percent of virus infections in 2004, according to a six-month\virus roundup published Wednesday by antivirus company Sophos."\\"The 18-year-old Jaschan was taken into custody in Germany in May by police who\said he had admitted programming both the Netsky and Sasser worms, something\experts at Microsoft confirmed. (A Microsoft antivirus reward program led to the\teenager''s arrest.) During the five months preceding Jaschan''s capture, there\were at least 25 variants of Netsky and one of the port-scanning network worm\Sasser."\\"Graham Cluley, senior technology consultant at Sophos, said it was staggeri ...\\',
  '\\FOAF/LOAF  and bloom filters have a lot of interesting properties for social\network and whitelist distribution.\\I think we can go one level higher though and include GPG/OpenPGP key\fingerpring distribution in the FOAF file for simple web-of-trust based key\distribution.\\What if we used FOAF and included the PGP key fingerprint(s) for identities?\This could mean a lot.  You include the PGP key fingerprints within the FOAF\file of your direct friends and then include a bloom filter of the PGP key\fingerprints of your entire whitelist (the source FOAF file would of course need\to be encrypted ).\\Your whitelist would be populated from the social network as your client\discovered new identit ...\\',
  'Wiltshire Police warns about "phishing" after its fraud squad chief was targeted.',
  'In its first two years, the UK''s dedicated card fraud unit, has recovered 36,000 stolen cards and 171 arrests - and estimates it saved 65m.',
  ' LOS ANGELES (Reuters) - A group of technology companies  including Texas Instruments Inc. &lt',
  ' LOS ANGELES (Reuters) - Apple Computer Inc.&lt',
  ' AMSTERDAM (Reuters) - Free Record Shop, a Dutch music  retail chain, beat Apple Computer Inc. to market on Tuesday  with the launch of a new download service in Europe''s latest  battleground for digital song services.',
  'A giant 100km colony of ants  which has been discovered in Melbourne, Australia, could threaten local insect species.',
  'Dolphin groups, or "pods", rely on socialites to keep them from collapsing, scientists claim.',
  'Tyrannosaurus rex achieved its massive size due to an enormous growth spurt during its adolescent years.',
  'Jet Propulsion Lab -- Scientists have discovered irregular lumps beneath the icy surface of Jupiter''s largest moon, Ganymede. These irregular masses may be rock formations, supported by Ganymede''s icy shell for billions of years...',
  'European Space Agency -- ESAs Mars Express has relayed pictures from one of NASA''s Mars rovers for the first time, as part of a set of interplanetary networking demonstrations.     The demonstrations pave the way for future Mars missions to draw on joint interplanetary networking capabilities...',
  'When did life begin? One evidential clue stems from the fossil records in Western Australia, although whether these layered sediments are biological or chemical has spawned a spirited debate. Oxford researcher, Nicola McLoughlin, describes some of the issues in contention.',
  'update Earnings per share rise compared with a year ago, but company misses analysts'' expectations by a long shot.',
  'By the end of the year, the computing giant plans to have its biggest headcount since 1991.',
  'Developers get early code for new operating system ''skin'' still being crafted.',
  'New technology applies electrical fuses to help identify and repair faults.',
  'Google has billed its IPO as a way for everyday people to get in on the process, denying Wall Street the usual stranglehold it''s had on IPOs. Public bidding, a minimum of just five shares, an open process with 28 underwriters - all this pointed to a new level of public participation. But this isn''t the case.',
  'By MICHAEL LIEDTKE     SAN FRANCISCO (AP) -- With its low prices and iconoclastic attitude, discount stock broker Charles Schwab Corp. (SCH) represented an annoying stone in Wall Street''s wing-tipped shoes for decades...',
  'Industry cyber security standards fail to reach some of the most vulnerable components of the power grid.\',
  'Michael Phelps won the gold medal in the 400 individual medley and set a world record in a time of 4 minutes 8.26 seconds.',
  'FOXBOROUGH -- Looking at his ridiculously developed upper body, with huge biceps and hardly an ounce of fat, it''s easy to see why Ty Law, arguably the best cornerback in football, chooses physical play over finesse. That''s not to imply that he''s lacking a finesse component, because he can shut down his side of the field much as Deion Sanders ...',
  'With the weeks dwindling until Jason Varitek enters free agency, the Red Sox continue to carefully monitor Kelly Shoppach , their catcher of the future, in his climb toward the majors. The Sox like most of what they have seen at Triple A Pawtucket from Shoppach, though it remains highly uncertain whether he can make the adjustments at the plate ...',
  'Just imagine what David Ortiz could do on a good night''s rest. Ortiz spent the night before last with his baby boy, D''Angelo, who is barely 1 month old. He had planned on attending the Red Sox'' Family Day at Fenway Park yesterday morning, but he had to sleep in. After all, Ortiz had a son at home, and he ...',
  'In  quot',
  'The Cleveland Indians pulled within one game of the AL Central lead by beating the Minnesota Twins, 7-1, Saturday night with home runs by Travis Hafner and Victor Martinez.',
  'Canadian Press - VANCOUVER (CP) - The sister of a man who died after a violent confrontation with police has demanded the city''s chief constable resign for defending the officer involved.',
  'NAJAF, Iraq - Explosions and gunfire rattled through the city of Najaf as U.S. troops in armored vehicles and tanks rolled back into the streets here Sunday, a day after the collapse of talks - and with them a temporary cease-fire - intended to end the fighting in this holy city...',
  'LOURDES, France - A frail Pope John Paul II, breathing heavily and gasping at times, celebrated an open-air Mass on Sunday for several hundred thousand pilgrims, many in wheelchairs, at a shrine to the Virgin Mary that is associated with miraculous cures.    At one point he said "help me" in Polish while struggling through his homily in French...',
  'Supporters and rivals warn of possible fraud',
  'AP - A 1994 law strengthened job protections for National Guard and Reserve troops called to active duty. Here are major provisions of the Uniformed Services Employment and Reemployment Rights Act (USERRA).',
  ' TEHRAN (Reuters) - A senior Iranian military official said  Sunday Israel and the United States would not dare attack Iran  since it could strike back anywhere in Israel with its latest  missiles, news agencies reported.',
  'KABUL, Afghanistan - Government troops intervened in Afghanistan''s latest outbreak of deadly fighting between warlords, flying from the capital to the far west on U.S. and NATO airplanes to retake an air base contested in the violence, officials said Sunday...',
  'AP - Randy Johnson took a four-hitter into the ninth inning to help the Arizona Diamondbacks end a nine-game losing streak Sunday, beating Steve Trachsel and the New York Mets 2-0.',
  'Reuters - Apparel retailers are hoping their\back-to-school fashions will make the grade among\style-conscious teens and young adults this fall, but it could\be a tough sell, with students and parents keeping a tighter\hold on their wallets.',
  'AP - If Hurricane Charley had struck three years ago, President Bush''s tour through the wreckage of this coastal city would have been just the sort of post-disaster visit that other presidents have made to the scenes of storms, earthquakes, floods and fires.',
  'FT.com - Shares in Sohu.com, a leading US-listed Chinese internet portal, fell more than 10 per cent on Friday after China''s biggest mobile phone network operator imposed a one-year suspension on its multimedia messaging services because of customers being sent spam.',
  'AP - Darin Erstad doubled in the go-ahead run in the eighth inning, lifting the Anaheim Angels to a 3-2 victory over the Detroit Tigers on Sunday. The win pulled Anaheim within a percentage point of Boston and Texas in the AL wild-card race.',
  'AP - Outfielder J.D. Drew missed the Atlanta Braves'' game against the St. Louis Cardinals on Sunday night with a sore right quadriceps.',
  ' CARACAS, Venezuela (Reuters) - Venezuelans voted in huge  numbers on Sunday in a historic referendum on whether to recall  left-wing President Hugo Chavez and electoral authorities  prolonged voting well into the night.',
  ' HONG KONG (Reuters) - Dell Inc. &lt',
  ' BEIJING (Reuters) - Beijing on Monday accused a  Chinese-American arrested for spying for Taiwan of building an  espionage network in the United States, and said he could go on  trial very soon.',
  'Another major, another disappointment for Tiger Woods, the No. 1 ranked player in the world who has not won a major championship since his triumph at the 2002 U.S. Open.',
  'AFP - A squadron of US Air Force F-15E fighters based in Alaska will fly to South Korea next month for temporary deployment aimed at enhancing US firepower on the Korean peninsula, US authorities said.',
  ' NEW YORK (Reuters) - Randy Johnson struck out 14 batters in  8 1/3 innings to help the Arizona Diamondbacks end a nine-game  losing streak with a 2-0 win over the host New York Mets in the  National League Sunday.',
  'AFP - A curfew in the capital of the Maldives was eased but parliament sessions were put off indefinitely and emergency rule continued following last week''s riots, officials and residents said.',
  'TheDeal.com - The U.K. mobile giant wants to find a way to disentagle the Czech wireless and fixed-line businesses.',
  ' LONDON (Reuters) - The dollar dipped to a four-week low  against the euro on Monday before rising slightly on  profit-taking, but steep oil prices and weak U.S. data  continued to fan worries about the health of the world''s  largest economy.',
  'As Michael Kaleko kept running into people who were getting older and having more vision problems, he realized he could do something about it.',
  'AFP - India''s Tata Iron and Steel Company Ltd. took a strategic step to expand its Asian footprint with the announcement it will buy the Asia-Pacific steel operations of Singapore''s NatSteel Ltd.',
  'BAGHDAD, Iraq - Delegates at Iraq''s National Conference called on radical Shiite cleric Muqtada al-Sadr to abandon his uprising against U.S. and Iraqi troops and pull his fighters out of a holy shrine in Najaf...',
  ' NEW YORK (Reuters) - U.S. Treasury debt prices slipped on  Monday, though traders characterized the move as profit-taking  rather than any fundamental change in sentiment.',
  ' NEW YORK (Reuters) - The dollar extended gains against the  euro on Monday after a report on flows into U.S. assets showed  enough of a rise in foreign investments to offset the current  account gap for the month.',
  ' MILWAUKEE (Sports Network) - U.S. Ryder Cup captain Hal  Sutton finalized his team on Monday when he announced the  selections of Jay Haas and Stewart Cink as his captain''s picks.',
  'Jay Haas joined Stewart Cink as the two captain''s picks for a U.S. team that will try to regain the cup from Europe next month.',
  'AP - American Natalie Coughlin won Olympic gold in the 100-meter backstroke Monday night. Coughlin, the only woman ever to swim under 1 minute in the event, finished first in 1 minute, 0.37 seconds. Kirsty Coventry of Zimbabwe, who swims at Auburn University in Alabama, earned the silver in 1:00.50. Laure Manaudou of France took bronze in 1:00.88.',
  'NewsFactor - Oracle (Nasdaq: ORCL) has revamped its sales-side CRM applications in version 11i.10 of its sales, marketing, partner relationship management and e-commerce application.',
  'AFP - The United Nations launched an appeal here for 210 million dollars to help flood victims facing "grave" food shortages after two-thirds of Bangladesh was submerged, destroying crops and killing more than 700 people.',
  'Government in South Indian state of Kerala sets up wireless kiosks as part of initiative to bridge digital divide.',
  'PUNTA GORDA, Fla. - Urban rescue teams, insurance adjusters and National Guard troops scattered across Florida Monday to help victims of Hurricane Charley and deliver water and other supplies to thousands of people left homeless...',
  'SANTA MARIA, Calif. - Fans of Michael Jackson erupted in cheers Monday as the pop star emerged from a double-decker tour bus and went into court for a showdown with the prosecutor who has pursued him for years on child molestation charges...',
  'AP - The Charlotte Bobcats traded center Predrag Drobnjak to the Atlanta Hawks on Monday for a second round pick in the 2005 NBA draft.',
  'Canadian Press - LANGLEY, B.C. (CP) - Police have arrested a man in the kidnapping and sexual assault of an 11-year-old girl that frightened this suburban Vancouver community last week.',
  'Red Flag Software Co., the company behind China''s leading Linux client distribution, plans to focus more on its server operating system and enterprise customers, the company''s acting president said.',
  'AOL Properties Sign Girafa For Thumbnail Search Images\\Girafa.com Inc. announced today that the CompuServe, Netscape, AIM and ICQ properties of America Online, Inc., have signed an agreement with Girafa to use Girafa''s thumbnail search images as an integrated part of their search results.\\Using Girafa''s thumbnail search service, search users can ...',
  'AP - NASA''s Cassini spacecraft has spied two new little moons around satellite-rich Saturn, the space agency said Monday.',
  'An industrial city northwest of Moscow struggles as AIDS hits a broader population.',
  'AP - A Nobel laureate in medicine said Monday the Bush administration''s limits on funding for embryonic stem cell research effectively have stopped the clock on American scientists'' efforts to develop treatments for a host of chronic, debilitating diseases.',
  'AP - Prosecutors suffered another setback Monday in the Kobe Bryant sexual assault case, losing a last-ditch attempt to keep the NBA star''s lawyers from telling jurors about the alleged victim''s sex life.',
  'Reuters - China has said no date has been set for\working-level talks on the North Korean nuclear crisis and gave\no indication that the meeting has been canceled, Australian\Foreign Minister Alexander Downer said on Tuesday.',
  'The Redskins expect huge things from 300-pound Cornelius Griffin, who was signed to aid the team''s weakest unit - the defensive line.',
  'AFP - The last surviving American defector to communist North Korea wants to tell his story to put a human face on the Stalinist state which he believes is unfairly vilified abroad, British film-makers said.',
  'Richard Faulds and Stephen Parry are going for gold for Great Britain on day four in Athens.',
  'Reuters - Democratic challenger John Kerry\has a commanding lead over President Bush in California of 54\percent to 38 percent among likely voters, a poll released on\Tuesday found.',
  ' ATHENS (Reuters) - At the beach volleyball, the 2004  Olympics is a sell-out, foot-stomping success.',
  ' LONDON (Reuters) - The dollar held steady near this week''s  four-week low against the euro on Tuesday with investors  awaiting a German investor confidence survey and U.S. consumer  inflation numbers to shed light on the direction.',
  'SAN FRANCISCO -- In the latest of a series of product delays, Intel Corp. has postponed the launch of a video display chip it had previously planned to introduce by year end, putting off a showdown with Texas Instruments Inc. in the fast-growing market for high-definition television displays.',
  'CARACAS -- Venezuelans voted resoundingly to keep firebrand populist Hugo Chavez as their president in a victory that drew noisy reactions yesterday from both sides in the streets. International observers certified the results as clean and accurate.',
  'AFP - Hong Kong democrats accused China of jailing one of their members on trumped-up prostitution charges in a bid to disgrace a political movement Beijing has been feuding with for seven years.',
  'AP - Former chess champion Bobby Fischer''s announcement thathe is engaged to a Japanese woman could win him sympathy among Japanese officials and help him avoid deportation to the United States, his fiancee and one of his supporters said Tuesday.',
  'ATHENS, Greece - Top American sprinters Jason Lezak and Ian Crocker missed the cut in the Olympic 100-meter freestyle preliminaries Tuesday, a stunning blow for a country that had always done well in the event.    Pieter van den Hoogenband of the Netherlands and Australian Ian Thorpe advanced to the evening semifinal a day after dueling teenager Michael Phelps in the 200 freestyle, won by Thorpe...',
  'A proposal backed by a coalition of telephone carriers would cut billions of dollars in fees owed by long-distance companies to regional phone giants but would allow the regional companies to make up some of the difference by raising monthly phone bills for millions of consumers. &lt',
  'Keep an eye on your credit card issuers -- they may be about to raise your rates.',
  'In another product postponement, semiconductor giant Intel Corp. said it won''t be offering a chip for projection TVs by the end of 2004 as it had announced earlier this year.',
  'NEW YORK (CNN/Money) - Money managers are growing more pessimistic about the economy, corporate profits and US stock market returns, according to a monthly survey by Merrill Lynch released Tuesday. ',
  'Olympic champion Kostas Kederis today left hospital ahead of his date with IOC inquisitors claiming his innocence and vowing:  quot',
  'NEWCASTLE, England (AP) - Striker Emile Heskey has pulled out of the England squad ahead of Wednesday #39',
  ' NEW YORK (Reuters) - Staples Inc. &lt',
  'AGHDAD, Iraq, Aug. 17  A delegation of Iraqis was delayed for security reasons today but still intended to visit Najaf to try to convince a rebellious Shiite cleric and his militia to evacuate a shrine in the holy city and end ...',
  ' WASHINGTON (Reuters) - U.S. consumer prices dropped in July  for the first time in eight months as a sharp run up in energy  costs reversed, the government said in a report that suggested  a slow rate of interest rate hikes is likely.',
  'An Indian army major shot his way to his country #39',
  'Rising fuel prices, a bugbear for most of the retail sector, are helping Home Depot (HD:NYSE - news - research), the remodeling giant that reported a surge in second-quarter earnings Tuesday and guided the rest of the year higher. ',
  'Charly Travers offers advice on withstanding the volatility of the biotech sector.',
  'Just what Alexander Downer was thinking when he declared on radio last Friday that  quot',
  ' ATHENS (Reuters) - World 100 meters champion Torri Edwards  will miss the Athens Olympics after her appeal against a  two-year drugs ban was dismissed on Tuesday, a source told  Reuters.',
  'NEW YORK - Stocks rose for a second straight session Tuesday as a drop in consumer prices allowed investors to put aside worries about inflation, at least for the short term.    With gasoline prices falling to eight-month lows, the Consumer Price Index registered a small drop in July, giving consumers a respite from soaring energy prices...',
  'Ilias Iliadis of Greece thrilled the home crowd Tuesday, beating Roman Gontyuk of Ukraine to win the gold medal in the 81-kilogram class. ',
  'AFP - Sudan will take the lead in restoring order to its rebellious Darfur region but needs the support of African peacekeepers and humanitarian aid, Foreign Minister Mustafa Osman Ismail said.',
  'The battle over home entertainment networking is heating up as a coalition proposes yet another standard for the IEEE #39',
  'Web giant Yahoo! is gambling that price cuts on its domain name registration and Web hosting products will make it more competitive with discounters in the space -- which means that small businesses looking to move online get a sweeter deal through ...',
  'IBM said Tuesday it has acquired a pair of Danish IT services firms as part of its effort to broaden its presence in Scandinavia. As a result of the moves, IBM will add about 3,700 IT staffers to its global head count. Financial terms of ...',
  'Motorola plans to sell mobile phone network equipment that uses Linux-based code, a step forward in network gear makers #39',
  'Microsoft will delay the release of its SP2 update for another week to fix software glitches. But not everyone is quite so eager to install the SP2 update for Windows XP. In fact, many companies have demanded the ability to prevent their ...',
  'Reuters - Two new moons were spotted around\Saturn by the Cassini space probe, raising the total to 33\moons for the ringed planet, NASA said on Monday.',
  'There are remarkable similarities between the 2004 Ohio State Buckeyes and those that won the national championship just two years ago. ',
  'The new IBM Power5 eServer i5 550 also features higher performance and new virtualization capabilities that allow it to run multiple operating systems at once on separate partitions.',
  'Newsday #146',
  'ATHENS -- The mistakes were so minor. Carly Patterson #39',
  'The price of oil has continued its sharp rise overnight, closing at a record high. The main contract in New York, light sweet crude for delivery next month, has closed at a record \$US46.75 a barrel - up 70 cents on yesterday #39',
  ' quot',
  'It wasn #39',
  'LONDON, AUGUST 17: Britain charged eight terror suspects on Tuesday with conspiracy to commit murder and said one had plans that could be used in striking US buildings that were the focus of security scares this month. ',
  'NewsFactor - IBM (NYSE: IBM) has -- again -- sought to have the pending legal claims by The SCO Group dismissed. According to a motion it filed in a U.S. district court, IBM argues that SCO has no evidence to support its claims that it appropriated confidential source code from Unix System V and placed it in Linux.',
  'NEW YORK - The newly released traffic crash fatality data have something for everyone in the debate about the safety of sport utility vehicles. ',
  'A CANADIAN husband #39',
  '-- The United States men #39',
  'BAGHDAD, Iraq - A national political conference #39',
  'TBILISI, Georgia  Georgian President Mikhail Saakashvili appealed to world leaders Tuesday to convene an international conference on the conflict in breakaway South Ossetia, where daily exchanges of gunfire threaten to spark ...',
  'AFP - Georgian and South Ossetian forces overnight accused each other of trying to storm the other side''s positions in Georgia''s breakaway region of South Ossetia, as four Georgian soldiers were reported to be wounded.',
  'BOSTON -- It was another busy day on the medical front for the Red Sox, as a series of roster moves were announced prior to Tuesday night #39',
  'AP - John Kerry, Bob Kerrey. It''s easy to get confused.',
  'AP - William H. Harlan, the retired University of Florida swimming coach who led the Gators to eight conference titles, died Tuesday, school officials said. He was 83.',
  'THENS, Aug. 17 - So Michael Phelps is not going to match the seven gold medals won by Mark Spitz. And it is too early to tell if he will match Aleksandr Dityatin, the Soviet gymnast who won eight total medals in 1980. But those were not the ...',
  'Three-year-old Victoria, from St Petersburg, has been living at the Schrders #39',
  'AP - Orlando Cabrera hit a run-scoring double off the Green Monster in the ninth inning on reliever Justin Speier''s second pitch of the game, giving the Boston Red Sox a 5-4 win over the Toronto Blue Jays on Tuesday night.',
  'Israel announced plans for 1,000 houses in the West Bank yesterday, accelerating the expansion of the settlements. ',
  'AP - At least one member of the top-ranked Southern California football team is under investigation for sexual assault, the Los Angeles Police Department said Tuesday.',
  'President Bush, in Pennsylvania, said that opponents of a missile defense system were putting the nation''s security at risk.',
  'BEIJING (Reuters) - China breathed a measured sigh of relief after the skills of its basketball giant Yao Ming dwarfed New Zealand to sweep his team nearer to their goal of reaching the Athens Olympics semi-finals. ',
  'A leaked Israeli plan to build 1,000 new Jewish settler homes in the West Bank yesterday sent Bush administration officials scrambling for a response in the sensitive period before November #39',
  'LONDON - British police charged eight terrorist suspects yesterday with conspiring to commit murder and use radioactive materials, toxic gases, chemicals or explosives to cause  quot',
  'Islamic group #39',
  ' TOKYO (Reuters) - The dollar moved in tight ranges on  Wednesday as most investors shrugged off lower-than-expected  U.S. inflation data and stuck to the view the U.S. Federal  Reserve would continue raising rates.',
  'Right-hander Matt Morris threw seven solid innings, but the Cardinals needed a bases-loaded walk to second baseman Tony Womack and a grand slam from new right fielder Larry Walker to key a six-run eighth inning for a ...',
  'ATHENS (Reuters) - Greek sprinters Costas Kenteris and Katerina Thanou have arrived at an Athens hotel for an International Olympic Committee (IOC) hearing into their missed doped tests, a saga that has shamed and angered the Olympic host ...',
  'BOSTON -- The Toronto Blue Jays have had worse hitting games this season against lesser pitchers than Pedro Martinez. '
)

$aVals = @(
  'Label',
  3,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  4,
  2,
  2,
  2,
  2,
  2,
  2,
  1,
  1,
  1,
  1,
  1,
  1,
  1,
  2,
  3,
  1,
  4,
  2,
  2,
  1,
  4,
  1,
  2,
  1,
  2,
  1,
  4,
  3,
  4,
  1,
  1,
  3,
  3,
  2,
  2,
  2,
  4,
  1,
  4,
  1,
  1,
  2,
  1,
  4,
  4,
  4,
  1,
  4,
  2,
  1,
  2,
  1,
  1,
  1,
  2,
  3,
  3,
  1,
  1,
  1,
  1,
  4,
  3,
  4,
  3,
  2,
  2,
  3,
  1,
  3,
  2,
  3,
  3,
  1,
  2,
  1,
  2,
  1,
  4,
  3,
  4,
  4,
  4,
  4,
  2,
  4,
  4,
  2,
  1,
  2,
  4,
  1,
  4,
  3,
  2,
  2,
  1,
  1,
  1,
  2,
  1,
  2,
  2,
  1,
  2,
  1,
  2,
  1,
  2,
  1,
  1,
  1,
  3,
  2,
  2,
  2
)

for ($i = 0; $i -lt $bText.Length; $i++) {
  $r = $i + 1
  $ws.Cells.Item($r, 1).Value = $aVals[$i]
  $ws.Cells.Item($r, 2).Value = $bText[$i]
}

$ws.Range("A144:B151").EntireRow.Delete()

[void]$ws.Range("C1:C1048576").Select()